$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 12 with 2021 data, following the pattern/formatting of row 11.
$ws.Range("A11:Y11").Copy()
$ws.Range("A12:Y12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A12").Value = "2021年"
$ws.Range("D12").Value = 2
$ws.Range("F12").Value = 686
$ws.Range("H12").Value = 315
$ws.Range("K12").Value = 82
$ws.Range("L12").Value = 76
$ws.Range("M12").Value = 2
$ws.Range("O12").Value = 3
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 137
$ws.Range("T12").Value = 40
$ws.Range("W12").Value = 25
